# fix(publipostage): Correct status name
#
# Renames the "bleu" status label to "noir" and corrects the wording of the
# statut_name values so they consistently read "... postés ou publiés ..."
# instead of "... et / ou publication posté ...".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# statut_label (column B): "bleu" -> "noir"
$ws.Range("B4").Value = "noir"

# statut_name (column C): correct wording for each status group
$ws.Range("C2").Value = "résultat postés ou publiés"
$ws.Range("C3").Value = "résultat postés ou publiés"
$ws.Range("C4").Value = "pas de résultat postés ni publiés"
$ws.Range("C5").Value = "résultat postés ou publiés"
$ws.Range("C6").Value = "résultat postés ou publiés"
$ws.Range("C7").Value = "résultat postés ou publiés dans les 36 mois"
$ws.Range("C8").Value = "résultat postés ou publiés dans les 12 mois"
$ws.Range("C9").Value = "résultat postés ou publiés dans les 12 mois"
